# Update cryptocurrency price/volume figures (and restore the Algorand /
# InternetComputer row ordering) per the "cryptos list" refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.059.42"
$ws.Range("E2").Value = "  -2.79%  "
$ws.Range("D3").Value = "1.730.68"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.31"
$ws.Range("E5").Value = "  -5.12%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4889"
$ws.Range("E7").Value = "  +6.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3515"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.62"
$ws.Range("E9").Value = "  +4.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07290"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.050"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.05"
$ws.Range("E13").Value = "  -2.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.898"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "1.730.48"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.903"
$ws.Range("E16").Value = "  -3.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.41"
$ws.Range("E17").Value = "  -4.52%  "
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06407"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.64"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.700"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").Value = "27.104.91"
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.93"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.078"
$ws.Range("E25").Value = "  -3.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.95"
$ws.Range("E26").Value = "  -4.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.04"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "1.929.10"
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.092"
$ws.Range("E29").Value = "  -2.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.53"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.056"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09335"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.634"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.403"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06004"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02192"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.444"
$ws.Range("E37").Value = "  +6.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.02"
$ws.Range("E38").Value = "  -5.71%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.790"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2001"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6024"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.100"
$ws.Range("E43").Value = "  -6.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.513"
$ws.Range("E44").Value = "  -3.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.85"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.583"
$ws.Range("E46").Value = "  -3.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5666"
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.07"
$ws.Range("E48").Value = "  -3.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.859"
$ws.Range("E49").Value = "  -3.19%  "
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06643"
$ws.Range("E51").Value = "  -2.27%  "
